# Repull data, push all data, mean calculation
# Update the dSF column (F) values for specific rows to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F7").Value = -5
$ws.Range("F9").Value = -2
$ws.Range("F12").Value = -2
$ws.Range("F15").Value = 4
$ws.Range("F16").Value = 1
$ws.Range("F18").Value = -5
$ws.Range("F21").Value = -2
$ws.Range("F22").Value = 4
$ws.Range("F23").Value = -4
$ws.Range("F25").Value = -3
$ws.Range("F29").Value = -8
